$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph by scanning the Paragraphs
# collection (robust against any earlier content changes/shifts).
$jupiterIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Ver no Jupiter") {
        $jupiterIndex = $i
        break
    }
}

# The copyright/footer paragraph is the one right after it, and the
# blank paragraph right before it should be removed too.
$startIndex = $jupiterIndex - 1
$endIndex = $jupiterIndex + 1

$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($endIndex)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
